{"js": "// Replace each distinct old value with its corresponding new value(s), in document order.\n// A few old values occur more than once in the document and map to different new values\n// depending on which occurrence it is, so each entry carries an ordered list of\n// replacement texts (one per occurrence, in document order).\nconst replacements = [\n  [\"2025-01-15 Wednesday\", [\"2025-01-16 Thursday\"]],\n  [\"32\u00f76=5, 2\", [\"37\u00f74=9, 1\"]],\n  [\"71\u00f77=10, 1\", [\"23\u00f74=5, 3\"]],\n  [\"61\u00f78=7, 5\", [\"87\u00f77=12, 3\"]],\n  [\"25\u00f75=5, 0\", [\"30\u00f78=3, 6\"]],\n  [\"70\u00f72=35, 0\", [\"69\u00f78=8, 5\"]],\n  [\"58\u00f74=14, 2\", [\"46\u00f75=9, 1\", \"48\u00f79=5, 3\"]],\n  [\"49\u00f76=8, 1\", [\"79\u00f73=26, 1\"]],\n  [\"24\u00f77=3, 3\", [\"50\u00f79=5, 5\"]],\n  [\"92\u00f78=11, 4\", [\"11\u00f75=2, 1\"]],\n  [\"60\u00f76=10, 0\", [\"76\u00f73=25, 1\", \"63\u00f73=21, 0\"]],\n  [\"94\u00f75=18, 4\", [\"24\u00f73=8, 0\"]],\n  [\"57\u00f73=19, 0\", [\"42\u00f79=4, 6\"]],\n  [\"80\u00f72=40, 0\", [\"53\u00f72=26, 1\"]],\n  [\"37\u00f75=7, 2\", [\"90\u00f73=30, 0\"]],\n  [\"20\u00f77=2, 6\", [\"36\u00f75=7, 1\"]],\n  [\"77\u00f73=25, 2\", [\"91\u00f78=11, 3\"]],\n  [\"30\u00f79=3, 3\", [\"12\u00f72=6, 0\"]],\n  [\"57\u00f78=7, 1\", [\"61\u00f78=7, 5\"]],\n  [\"74\u00f73=24, 2\", [\"29\u00f79=3, 2\"]],\n  [\"31\u00f75=6, 1\", [\"75\u00f76=12, 3\"]],\n  [\"55\u00f74=13, 3\", [\"10\u00f79=1, 1\"]],\n  [\"62\u00f76=10, 2\", [\"32\u00f75=6, 2\"]],\n  [\"31\u00f73=10, 1\", [\"72\u00f76=12, 0\"]],\n];\n\nfor (const [oldText, newTexts] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== newTexts.length) {\n    throw new Error(\n      `Expected ${newTexts.length} occurrence(s) of \"${oldText}\" but found ${results.items.length}`\n    );\n  }\n\n  // results.items[] is returned in document order, so pairing index-for-index with\n  // newTexts[] correctly disambiguates repeated old values (e.g. the two cells that\n  // originally read \"58\u00f74=14, 2\" become different values).\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each distinct old value with its corresponding new value(s), in document order.\n# Some old values occur more than once in the document and map to different new values\n# depending on occurrence order, so each entry carries an ordered list of replacement texts.\n$replacements = @(\n  @{ Old = \"2025-01-15 Wednesday\"; New = @(\"2025-01-16 Thursday\") },\n  @{ Old = \"32\u00f76=5, 2\"; New = @(\"37\u00f74=9, 1\") },\n  @{ Old = \"71\u00f77=10, 1\"; New = @(\"23\u00f74=5, 3\") },\n  @{ Old = \"61\u00f78=7, 5\"; New = @(\"87\u00f77=12, 3\") },\n  @{ Old = \"25\u00f75=5, 0\"; New = @(\"30\u00f78=3, 6\") },\n  @{ Old = \"70\u00f72=35, 0\"; New = @(\"69\u00f78=8, 5\") },\n  @{ Old = \"58\u00f74=14, 2\"; New = @(\"46\u00f75=9, 1\", \"48\u00f79=5, 3\") },\n  @{ Old = \"49\u00f76=8, 1\"; New = @(\"79\u00f73=26, 1\") },\n  @{ Old = \"24\u00f77=3, 3\"; New = @(\"50\u00f79=5, 5\") },\n  @{ Old = \"92\u00f78=11, 4\"; New = @(\"11\u00f75=2, 1\") },\n  @{ Old = \"60\u00f76=10, 0\"; New = @(\"76\u00f73=25, 1\", \"63\u00f73=21, 0\") },\n  @{ Old = \"94\u00f75=18, 4\"; New = @(\"24\u00f73=8, 0\") },\n  @{ Old = \"57\u00f73=19, 0\"; New = @(\"42\u00f79=4, 6\") },\n  @{ Old = \"80\u00f72=40, 0\"; New = @(\"53\u00f72=26, 1\") },\n  @{ Old = \"37\u00f75=7, 2\"; New = @(\"90\u00f73=30, 0\") },\n  @{ Old = \"20\u00f77=2, 6\"; New = @(\"36\u00f75=7, 1\") },\n  @{ Old = \"77\u00f73=25, 2\"; New = @(\"91\u00f78=11, 3\") },\n  @{ Old = \"30\u00f79=3, 3\"; New = @(\"12\u00f72=6, 0\") },\n  @{ Old = \"57\u00f78=7, 1\"; New = @(\"61\u00f78=7, 5\") },\n  @{ Old = \"74\u00f73=24, 2\"; New = @(\"29\u00f79=3, 2\") },\n  @{ Old = \"31\u00f75=6, 1\"; New = @(\"75\u00f76=12, 3\") },\n  @{ Old = \"55\u00f74=13, 3\"; New = @(\"10\u00f79=1, 1\") },\n  @{ Old = \"62\u00f76=10, 2\"; New = @(\"32\u00f75=6, 2\") },\n  @{ Old = \"31\u00f73=10, 1\"; New = @(\"72\u00f76=12, 0\") }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n  # Walk occurrences left-to-right by re-scoping the search range to start right after the\n  # previous match. This lets repeated Old values (e.g. \"58\u00f74=14, 2\") be replaced with\n  # different New values depending on which occurrence in the document it is, instead of\n  # Find.Execute's wdReplaceAll stamping every occurrence with the same text.\n  $searchStart = 0\n  foreach ($newText in $pair.New) {\n    $r = $d.Range($searchStart, $d.Content.End)\n    $r.Find.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap(1=wdFindStop), Format, ReplaceWith,\n    #          Replace(1=wdReplaceOne))\n    $found = $r.Find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n    if (-not $found) {\n      throw \"Could not find occurrence of $($pair.Old)\"\n    }\n    $searchStart = $r.End\n  }\n}"}
